$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.954.75'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.640.32'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.54'
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("D12").Value = '1.872.97'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '1.633.04'
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.573'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.91'
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("D17").Value = '27.953.89'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.53'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").Value = '0.0₃0724'
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.37'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.18'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.98'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.68'
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D34").Value = '1.412.22'
$ws.Range("E34").Value = '  -4.33%  '
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0169'
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.884'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.905'
$ws.Range("E40").Value = '  -4.75%  '
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.88'
$ws.Range("E43").Value = '  +6.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.53'
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("E45").Value = '  +3.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '1.781.69'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.84'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("E51").Value = '  -1.53%  '
